$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "41.505.08"
Set-TextValue "E2" "  -0.55%  "

Set-TextValue "D3" "2.196.60"
Set-TextValue "E3" "  -2.70%  "

Set-TextValue "E4" "  +0.00%  "

Set-TextValue "D5" "229.09"
Set-TextValue "E5" "  -1.75%  "

Set-TextValue "D6" "0.616"
Set-TextValue "E6" "  -4.31%  "

Set-TextValue "D7" "60.05"
Set-TextValue "E7" "  -6.71%  "

Set-TextValue "E8" "  -0.03%  "

Set-TextValue "D9" "0.400"
Set-TextValue "E9" "  -3.64%  "

Set-TextValue "D10" "56.72"
Set-TextValue "E10" "  -5.89%  "

Set-TextValue "D11" "0.0884"
Set-TextValue "E11" "  -2.67%  "

Set-TextValue "E12" "  -2.07%  "

Set-TextValue "D13" "2.521.86"
Set-TextValue "E13" "  -2.58%  "

Set-TextValue "D14" "15.32"
Set-TextValue "E14" "  -5.72%  "

Set-TextValue "D15" "22.16"
Set-TextValue "E15" "  -2.44%  "

Set-TextValue "D16" "5.63"
Set-TextValue "E16" "  -1.21%  "

Set-TextValue "D17" "0.790"
Set-TextValue "E17" "  -4.93%  "

Set-TextValue "D18" "2.196.16"
Set-TextValue "E18" "  -2.47%  "

Set-TextValue "D19" "41.409.94"
Set-TextValue "E19" "  -0.33%  "

Set-TextValue "D20" "71.83"
Set-TextValue "E20" "  -2.99%  "

Set-TextValue "D21" "0.0₃0895"
Set-TextValue "E21" "  -3.82%  "

Set-TextValue "D22" "6.02"
Set-TextValue "E22" "  -3.10%  "

Set-TextValue "D23" "241.72"
Set-TextValue "E23" "  -4.68%  "

Set-TextValue "D25" "2.35"
Set-TextValue "E25" "  -2.42%  "

Set-TextValue "D26" "2.29"
Set-TextValue "E26" "  -2.63%  "

Set-TextValue "D27" "9.56"
Set-TextValue "E27" "  -3.17%  "

Set-TextValue "D28" "168.65"
Set-TextValue "E28" "  -1.98%  "

Set-TextValue "E29" "  -7.12%  "

Set-TextValue "E30" "  -0.89%  "

Set-TextValue "D31" "19.65"
Set-TextValue "E31" "  -4.33%  "

Set-TextValue "D32" "2.56"
Set-TextValue "E32" "  -8.85%  "

Set-TextValue "E33" "  -4.38%  "

Set-TextValue "D34" "4.97"
Set-TextValue "E34" "  -3.28%  "

Set-TextValue "D35" "4.60"
Set-TextValue "E35" "  -3.76%  "

Set-TextValue "D36" "0.0645"
Set-TextValue "E36" "  +0.25%  "

Set-TextValue "B37" "LidoDAOToken"
Set-TextValue "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D37" "2.34"
Set-TextValue "E37" "  -5.21%  "

Set-TextValue "B38" "THORChain"
Set-TextValue "C38" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D38" "6.29"
Set-TextValue "E38" "  -9.21%  "

Set-TextValue "D39" "3.52"
Set-TextValue "E39" "  -8.44%  "

Set-TextValue "B40" "TerraClassic"
Set-TextValue "C40" "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue "D40" "0.000239"
Set-TextValue "E40" "  -7.63%  "

Set-TextValue "B41" "BinanceUSD"
Set-TextValue "C41" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D41" "0.999"
Set-TextValue "E41" "  -0.12%  "

Set-TextValue "D42" "0.0236"
Set-TextValue "E42" "  -3.25%  "

Set-TextValue "D43" "8.48"
Set-TextValue "E43" "  -3.85%  "

Set-TextValue "D44" "0.0949"
Set-TextValue "E44" "  -5.09%  "

Set-TextValue "B45" "TrustWalletToken"
Set-TextValue "C45" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.19"
Set-TextValue "E45" "  -3.55%  "

Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "96.65"
Set-TextValue "E46" "  -6.11%  "

Set-TextValue "D47" "1.459.17"
Set-TextValue "E47" "  -3.61%  "

Set-TextValue "E48" "  -16.01%  "

Set-TextValue "D49" "16.29"
Set-TextValue "E49" "  -8.07%  "

Set-TextValue "D50" "2.76"
Set-TextValue "E50" "  -1.35%  "

Set-TextValue "E51" "  -7.30%  "

